$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 - I11 existing entry (MVC EF learning comment) gets restyled
# (style 45 -> 44, font no longer the "charset 204" variant)
$ws.Range("I11").Style = "Normal"
$ws.Cells.Item(11, 9).Font.Name = "Calibri"

# Row 13 - new log entry: "Help"
$ws.Range("B13").Value = 6
$ws.Range("C13").Value = 43528
$ws.Range("D13").Value = 0.79861111111111116
$ws.Range("E13").Value = 0.88888888888888884
$ws.Range("F13").Value = 10
$ws.Range("G13").Value = 120
$ws.Range("H13").Value = "Help"
$ws.Range("I13").Value = "Helping my friend to solve problems with VS + MVC "
$ws.Rows.Item(13).RowHeight = 28.8

# Row 14 - new log entry: "Practise"
$ws.Range("B14").Value = 7
$ws.Range("C14").Value = 43525
$ws.Range("D14").Value = 0.33333333333333331
$ws.Range("E14").Value = 0.41666666666666669
$ws.Range("G14").Value = 120
$ws.Range("H14").Value = "Practise"
$ws.Range("I14").Value = "In class"

# Move the active selection like the author left it
$ws.Range("I17").Select()
